$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp cell (E25, "project:" footer row) before the row
# shift below changes which row index that footer lives on. The leading
# apostrophe forces text (matching how "16:30" was originally stored,
# i.e. with quotePrefix) rather than letting it parse as a time value.
$ws.Range("E25").Value = "'17:42"

# Remove the B1 / N-5 line item (row 23) entirely; rows below it shift
# up by one (old 24->23, 25->24, 26->25) and the ROW()-based numbering
# formulas + shared-string indices re-resolve automatically.
$ws.Rows(23).Delete()
